$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Relocate the existing header row (B2:H2) down to A5:G5 (shifted
#    one column left, three rows down) using Cut so the shared-string
#    indices / styles travel with the cells instead of being rebuilt.
# ------------------------------------------------------------------
$ws.Range("B2:H2").Cut($ws.Range("A5"))
$ws.Range("B2:H2").Clear()

# ------------------------------------------------------------------
# 2. Merge the cells that will hold the title / part number / description
#    before they carry any value or custom formatting.
# ------------------------------------------------------------------
$ws.Range("A1:G1").Merge()
$ws.Range("B2:D2").Merge()
$ws.Range("B3:D3").Merge()

# ------------------------------------------------------------------
# 3. Row 2 - "Part Number" label + value (written first so the new
#    shared strings land in the same order as the target file).
# ------------------------------------------------------------------
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").Value = "Part Number"
$partNumberValue = $ws.Range("B2:D2")
$partNumberValue.HorizontalAlignment = -4131   # xlLeft
$ws.Range("B2").Value = "19-ELC-1234"

# ------------------------------------------------------------------
# 4. Row 3 - "Description" label + value.
# ------------------------------------------------------------------
$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").Value = "Description"
$descriptionValue = $ws.Range("B3:D3")
$descriptionValue.HorizontalAlignment = -4131  # xlLeft
$ws.Range("B3").Value = "Mayonnaise Dispenser"

# ------------------------------------------------------------------
# 5. Row 1 - big bold centered title banner.
# ------------------------------------------------------------------
$title = $ws.Range("A1:G1")
$title.Font.Bold = $true
$title.Font.Size = 16
$title.HorizontalAlignment = -4108             # xlCenter
$ws.Range("A1").Value = "Longhorn Racing Electric BOM"
$ws.Rows(1).RowHeight = 21

# ------------------------------------------------------------------
# 6. Column widths.
# ------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 23.833333333333332
$ws.Columns("B").ColumnWidth = 5.5
$ws.Range("C1:G1").ColumnWidth = 19.666666666666668

# ------------------------------------------------------------------
# 7. Selection cursor parked below the table, matching the saved view.
# ------------------------------------------------------------------
[void]$ws.Range("A6").Select()
